$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "CubeA-HW50.xpc" to "CubeA"
$ws.Name = "CubeA"

# Add a new row 16, mirroring the pattern of row 15 (Gaussian Quadrature scheme export)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = $ws.Range("B15").Text

$ws.Range("C16:P16").Value = 1
